# "version afinada de dashboard" - refresh the rolling seismic data window:
# drop the oldest (2025-12-07) rows and append the newest (2025-12-10, 2025-12-11) readings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old last row (row 8) first, then overwrite rows 2-7 with the refreshed data.
$ws.Rows.Item(8).Delete()

$data = @(
    @("2025-12-09", "00:13:21", "-6.01", "-79.91", "76", "3.6"),
    @("2025-12-09", "02:46:54", "-7.25", "-80.28", "24", "4"),
    @("2025-12-09", "04:34:19", "-15.98", "-74.89", "24", "4.4"),
    @("2025-12-09", "18:56:00", "-10.68", "-74.68", "16", "3.6"),
    @("2025-12-10", "09:23:37", "-12.03", "-77.47", "52", "3.6"),
    @("2025-12-11", "19:59:41", "-13.61", "-72.15", "9", "3.5")
)

$ws.Range("A2:F7").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $rowData[$j]
    }
}

$ws.Range("A2:F7").Style = "Normal"
